# Apply "new parameters for capacity/reserves" edit to parameter.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) New data rows 138-145 ("capacity and reserves" parameter group)
#    Write values first, then paste formatting (styles only) from an
#    existing row (114) that already uses the target style combination
#    (A/B/E/F = style 27, C = style 14, G = style 28, D/H/I/J = style 29).
# ---------------------------------------------------------------------

$rows = @(
  @{ r=138; a="capacity and reserves"; b="car_sc";        c="kareb.vf";              d=25;    e="percent";       f="medium"; g="Proportion of staircases; for the conversion from total area to living area."; h=25;    i=25;    j=25 },
  @{ r=139; a="capacity and reserves"; b="car_resi";      c="kareb.wohnant";         d=0;     e="percent";       f="high";   g="Residence portion ('slider'): -100% = minimum residential share according to the building regulation (BZO), 0% = real residential share, +100% = maximum residential share according to building regulation."; h=-25; i=0; j=25 },
  @{ r=140; a="capacity and reserves"; b="car_plot";      c="kareb.areal";           d=50;    e="percent";       f="high";   g="Proportion of plot construction ('slider'); 0% = without plot construction, 100% = with plot construction"; h=0; i=50; j=100 },
  @{ r=141; a="capacity and reserves"; b="car_uti_input"; c="no parameter in previous model"; d=85; e="percent";  f="high";   g="degree of utilization used to calculate the input data"; h=85; i=85; j=85 },
  @{ r=142; a="capacity and reserves"; b="car_uti";       c="kareb.ausbau";          d=85;    e="percent";       f="high";   g="degree of utilization; linear influence on the capacity"; h=75; i=85; j=90 },
  @{ r=143; a="capacity and reserves"; b="car_pp";        c="kareb.ina.prozentpunkte"; d=0;   e="percent points"; f="medium"; g="usage: less or more percentage points (parameter between -100 and +100)"; h=0; i=0; j=0 },
  @{ r=144; a="capacity and reserves"; b="car_y";         c="kareb.ina.jahr";        d=2044;  e="year";          f="low";    g="Reference year of the usage values"; h=2044; i=2044; j=2044 },
  @{ r=145; a="capacity and reserves"; b="car_lamda";     c="kareb.ina.lambda";      d=-0.04; e="per year";      f="high";   g="lambda value of an exponential function exp(lambda * time since start of scenario); proportion of utilization per year"; h=-0.04; i=-0.04; j=-0.04 }
)

foreach ($row in $rows) {
    $r = $row.r
    $ws.Range("A$r").Value = $row.a
    $ws.Range("B$r").Value = $row.b
    $ws.Range("C$r").Value = $row.c
    $ws.Range("D$r").Value = $row.d
    $ws.Range("E$r").Value = $row.e
    $ws.Range("F$r").Value = $row.f
    $ws.Range("G$r").Value = $row.g
    $ws.Range("H$r").Value = $row.h
    $ws.Range("I$r").Value = $row.i
    $ws.Range("J$r").Value = $row.j

    $ws.Range("A114:J114").Copy()
    $ws.Range("A${r}:J${r}").PasteSpecial(-4122)
}
$excel.CutCopyMode = 0

# Row heights that Excel auto-computed for the wrapped long descriptions
$ws.Rows.Item(139).RowHeight = 38.25
$ws.Rows.Item(140).RowHeight = 25.5
$ws.Rows.Item(145).RowHeight = 25.5

# Existing row 55 picks up an auto height too once the sheet reflows
$ws.Rows.Item(55).RowHeight = 25.5

# ---------------------------------------------------------------------
# 2) Column F narrows (it no longer needs to fit the long description
#    text that moved to column G for these new rows)
# ---------------------------------------------------------------------
$ws.Columns.Item(6).ColumnWidth = 19

# ---------------------------------------------------------------------
# 3) View state: zoom, frozen-pane anchor and active selection moved to
#    show the newly added rows.
# ---------------------------------------------------------------------
$ws.Application.ActiveWindow.Zoom = 90
$ws.Application.ActiveWindow.ScrollRow = 123
$ws.Range("G143").Select()
